$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028598544292766
$ws.Range("D2").Value = 1.032458595806573
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.038891713623375
$ws.Range("I2").Value = 1.034979182921329
$ws.Range("J2").Value = 1.033749972674907
$ws.Range("K2").Value = 1.035263767413752
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.041678455416404
$ws.Range("N2").Value = 1.035218016208896

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029392318614242
$ws.Range("D3").Value = 1.033033649030302
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.039949570567998
$ws.Range("I3").Value = 1.035151277759758
$ws.Range("J3").Value = 1.034185221895733
$ws.Range("K3").Value = 1.035648222598183
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.04254576501394
$ws.Range("N3").Value = 1.035653883533528

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.02990634035713
$ws.Range("D4").Value = 1.03340607743262
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.04063489982492
$ws.Range("I4").Value = 1.035261662921154
$ws.Range("J4").Value = 1.034466583496681
$ws.Range("K4").Value = 1.035896619642251
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.043107199568469
$ws.Range("N4").Value = 1.035935644700222

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030122528573879
$ws.Range("D5").Value = 1.033562723876627
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.040923208785074
$ws.Range("I5").Value = 1.035307835623805
$ws.Range("J5").Value = 1.034584801329908
$ws.Range("K5").Value = 1.036000955854026
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.043343280099048
$ws.Range("N5").Value = 1.036054030416332

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03015883296495
$ws.Range("D6").Value = 1.033589030004991
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.040971628624842
$ws.Range("I6").Value = 1.035315574530521
$ws.Range("J6").Value = 1.034604646689321
$ws.Range("K6").Value = 1.036018469060522
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.0433829221651
$ws.Range("N6").Value = 1.036073903958432

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029909228709228
$ws.Range("D7").Value = 1.03340817024737
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.040638751452307
$ws.Range("I7").Value = 1.035262280800193
$ws.Range("J7").Value = 1.034468163392095
$ws.Range("K7").Value = 1.0358980141424
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.043110353879269
$ws.Range("N7").Value = 1.035937226839269

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028866721010128
$ws.Range("D8").Value = 1.032652868350585
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.039249049899418
$ws.Range("I8").Value = 1.035037544020466
$ws.Range("J8").Value = 1.033897123238612
$ws.Range("K8").Value = 1.035393772066699
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.041971519123387
$ws.Range("N8").Value = 1.035365375743268

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027032793171095
$ws.Range("D9").Value = 1.031324531061299
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.036806580813004
$ws.Range("I9").Value = 1.034634113117263
$ws.Range("J9").Value = 1.032888827199471
$ws.Range("K9").Value = 1.034502436635573
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.039966529663436
$ws.Range("N9").Value = 1.034355647808099

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025812349815972
$ws.Range("D10").Value = 1.030440815796689
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.035182603204893
$ws.Range("I10").Value = 1.034360208178974
$ws.Range("J10").Value = 1.032215313422033
$ws.Range("K10").Value = 1.033906394978934
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.038631124786036
$ws.Range("N10").Value = 1.033681177563843

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025284418769294
$ws.Range("D11").Value = 1.030058614273536
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.03448044281799
$ws.Range("I11").Value = 1.034240437342194
$ws.Range("J11").Value = 1.03192337529178
$ws.Range("K11").Value = 1.033647884190975
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.038053188781011
$ws.Range("N11").Value = 1.03338882484796

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025088402615006
$ws.Range("D12").Value = 1.029916717271119
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.034219785208494
$ws.Range("I12").Value = 1.034195774183815
$ws.Range("J12").Value = 1.031814892098902
$ws.Range("K12").Value = 1.033551799589083
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.037838564026577
$ws.Range("N12").Value = 1.033280187596505

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025130445059065
$ws.Range("D13").Value = 1.029947151481872
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.034275690076039
$ws.Range("I13").Value = 1.034205362489982
$ws.Range("J13").Value = 1.031838164103418
$ws.Range("K13").Value = 1.033572412859704
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.03788459965077
$ws.Range("N13").Value = 1.033303492649937

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025268214353125
$ws.Range("D14").Value = 1.030046883592752
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.034458893590763
$ws.Range("I14").Value = 1.034236749036941
$ws.Range("J14").Value = 1.031914408936436
$ws.Range("K14").Value = 1.033639943070365
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.038035446880046
$ws.Range("N14").Value = 1.033379845759362

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025353109331812
$ws.Range("D15").Value = 1.030108341116363
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.03457179202352
$ws.Range("I15").Value = 1.034256064160683
$ws.Range("J15").Value = 1.031961380045511
$ws.Range("K15").Value = 1.033681542452843
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.038128395012544
$ws.Range("N15").Value = 1.033426883572801

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025847398131223
$ws.Range("D16").Value = 1.03046619093426
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.035229225108675
$ws.Range("I16").Value = 1.034368132411193
$ws.Range("J16").Value = 1.032234682114963
$ws.Range("K16").Value = 1.03392354271103
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.038669486955837
$ws.Range("N16").Value = 1.033700573762539

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02615759509587
$ws.Range("D17").Value = 1.03069078308725
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.03564189259892
$ws.Range("I17").Value = 1.034438117508443
$ws.Range("J17").Value = 1.032406037098585
$ws.Range("K17").Value = 1.034075231147703
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.039008981402812
$ws.Range("N17").Value = 1.033872172089892

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026338578677505
$ws.Range("D18").Value = 1.030821827492883
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.03588269413905
$ws.Range("I18").Value = 1.034478825879774
$ws.Range("J18").Value = 1.032505956332364
$ws.Range("K18").Value = 1.034163667774117
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.039207031957915
$ws.Range("N18").Value = 1.033972233220444

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026400297997477
$ws.Range("D19").Value = 1.030866517595734
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.035964818104292
$ws.Range("I19").Value = 1.034492687240388
$ws.Range("J19").Value = 1.032540021219765
$ws.Range("K19").Value = 1.03419381544799
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.039274567003508
$ws.Range("N19").Value = 1.034006346483893

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026124308617631
$ws.Range("D20").Value = 1.030666681953179
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.035597606976767
$ws.Range("I20").Value = 1.034430620428171
$ws.Range("J20").Value = 1.03238765534637
$ws.Range("K20").Value = 1.034058960613446
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.038972553836225
$ws.Range("N20").Value = 1.033853764233481

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025227642483616
$ws.Range("D21").Value = 1.030017513036047
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.034404940407837
$ws.Range("I21").Value = 1.034227511305673
$ws.Range("J21").Value = 1.03189195795008
$ws.Range("K21").Value = 1.033620058814752
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.037991024842624
$ws.Range("N21").Value = 1.033357362890031

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024664341933893
$ws.Range("D22").Value = 1.029609758597521
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.033655966733966
$ws.Range("I22").Value = 1.034098796692173
$ws.Range("J22").Value = 1.031580037106666
$ws.Range("K22").Value = 1.033343744563662
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.037374167558581
$ws.Range("N22").Value = 1.033044999083242

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024962913262929
$ws.Range("D23").Value = 1.029825878129166
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.034052925947253
$ws.Range("I23").Value = 1.034167126464584
$ws.Range("J23").Value = 1.031745416189164
$ws.Range("K23").Value = 1.033490257684487
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.037701149490821
$ws.Range("N23").Value = 1.033210613023006

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026139349187116
$ws.Range("D24").Value = 1.03067757208372
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.035617617435531
$ws.Range("I24").Value = 1.034434008384485
$ws.Range("J24").Value = 1.032395961359392
$ws.Range("K24").Value = 1.034066312693311
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.038989013795416
$ws.Range("N24").Value = 1.033862082041994

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.0275065305519
$ws.Range("D25").Value = 1.031667620135031
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.037437257736437
$ws.Range("I25").Value = 1.034739284898112
$ws.Range("J25").Value = 1.033149732468053
$ws.Range("K25").Value = 1.034733193198442
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.040484650458637
$ws.Range("N25").Value = 1.034616923592089
